$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "80÷4=20, 0"
$t.Cell(1, 2).Range.Text = "50÷2=25, 0"
$t.Cell(1, 3).Range.Text = "88÷5=17, 3"
$t.Cell(1, 4).Range.Text = "79÷9=8, 7"
$t.Cell(1, 5).Range.Text = "79÷3=26, 1"
$t.Cell(5, 1).Range.Text = "87÷5=17, 2"
$t.Cell(5, 2).Range.Text = "72÷8=9, 0"
$t.Cell(5, 3).Range.Text = "55÷3=18, 1"
$t.Cell(5, 4).Range.Text = "96÷4=24, 0"
$t.Cell(5, 5).Range.Text = "58÷2=29, 0"
$t.Cell(9, 1).Range.Text = "24÷2=12, 0"
$t.Cell(9, 2).Range.Text = "90÷8=11, 2"
$t.Cell(9, 3).Range.Text = "41÷4=10, 1"
$t.Cell(9, 4).Range.Text = "58÷6=9, 4"
$t.Cell(9, 5).Range.Text = "16÷2=8, 0"
$t.Cell(13, 1).Range.Text = "39÷6=6, 3"
$t.Cell(13, 2).Range.Text = "86÷7=12, 2"
$t.Cell(13, 3).Range.Text = "45÷8=5, 5"
$t.Cell(13, 4).Range.Text = "87÷7=12, 3"
$t.Cell(13, 5).Range.Text = "23÷5=4, 3"
$t.Cell(17, 1).Range.Text = "32÷3=10, 2"
$t.Cell(17, 2).Range.Text = "26÷4=6, 2"
$t.Cell(17, 3).Range.Text = "22÷3=7, 1"
$t.Cell(17, 4).Range.Text = "13÷3=4, 1"
$t.Cell(17, 5).Range.Text = "47÷4=11, 3"
